# Applies the crypto price / 1h-volume refresh described in the commit
# "Updated cryptos list on Sun May 12 05:34:11 UTC 2024 with GitHub Actions".
# All source values are plain text in the workbook (prices use "." as a
# thousands separator in some rows, e.g. "61.028.36"), so every write below
# that targets a plain decimal-looking price (e.g. "1.00", "6.70") first
# forces the cell's NumberFormat to Text ("@"). Otherwise Excel's normal
# type-inference on Range.Value would silently reinterpret it as the
# number 1 or 6.7 and drop the significant trailing zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.993.97'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '2.923.85'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '590.60'
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.95'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.95'
$ws.Range('E9').Value = '  +0.73%  '
$ws.Range('E10').Value = '  -0.55%  '
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '33.69'
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').Value = '3.407.56'
$ws.Range('E15').Value = '  +0.21%  '
$ws.Range('D16').Value = '60.927.68'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.925.58'
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.70'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '431.71'
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('E20').Value = '  -1.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.679'
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('E22').Value = '  -0.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '81.38'
$ws.Range('E23').Value = '  +1.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.92'
$ws.Range('E24').Value = '  +0.74%  '
$ws.Range('E25').Value = '  -0.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.89'
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.26'
$ws.Range('E28').Value = '  +4.63%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.01'
$ws.Range('E30').Value = '  -3.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.68'
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.109'
$ws.Range('E32').Value = '  +2.18%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').Value = '0.0₃0865'
$ws.Range('E34').Value = '  -1.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.01'
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.00'
$ws.Range('E37').Value = '  -1.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.98'
$ws.Range('E38').Value = '  -1.18%  '
$ws.Range('E39').Value = '  -5.09%  '
$ws.Range('E40').Value = '  -1.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '41.36'
$ws.Range('E41').Value = '  +0.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.282'
$ws.Range('E42').Value = '  -4.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '378.06'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').Value = '2.706.89'
$ws.Range('E44').Value = '  +0.97%  '
$ws.Range('E45').Value = '  -1.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '134.08'
$ws.Range('E47').Value = '  -0.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.81'
$ws.Range('E48').Value = '  -3.99%  '
$ws.Range('E49').Value = '  -0.59%  '
$ws.Range('E50').Value = '  -2.88%  '
$ws.Range('E51').Value = '  -0.71%  '
